$d = $word.ActiveDocument

# 1. Update the ASIC sentence in the mining section from Vietnamese to the
#    (untranslated) English text, matching the committed revision.
$old = "Khai thác mỏ SmartCash ngăn tập trung khai thác mỏ và kích thích phát triển mạng lưới. Mỗi máy tính có thể được sử dụng như một thiết bị khai thác mỏ trong khi vẫn cho phép máy tính đó được sử dụng cho các nhiệm vụ khác. ASIC vẫn chưa được tạo ra cho thuật toán khai thác Keccak và có thể an toàn để giả định rằng không có ASIC nào được tạo ra trong một thời gian dài, cho đến khi Smartcash đạt đến mức vốn hóa thị trường đáng kể."
$new = "Khai thác mỏ SmartCash ngăn tập trung khai thác mỏ và kích thích phát triển mạng lưới. Mỗi máy tính có thể được sử dụng như một thiết bị khai thác mỏ trong khi vẫn cho phép máy tính đó được sử dụng cho các nhiệm vụ khác. ASICs have yet to be created for the Keccak mining algorithm and it" + [char]0x2019 + "s probably safe to assume no ASICs will be created for quite some time."

$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

Write-Output "done"
